# Add 5 new participant rows (blg087, blg088, blg092, blg095, blg097) to Sheet1,
# each with Age=5 and Range=5, following the same pattern/format as the
# preceding rows (e.g. row 97).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @("blg087", "blg088", "blg092", "blg095", "blg097")
$startRow = 98

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $startRow + $i

    # Copy the formatting (styles) from the row above (the last existing
    # data row) so the new rows look consistent with the rest of the table.
    $ws.Range("A" + ($row - 1) + ":C" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":C" + $row).PasteSpecial(-4122) # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = 5
    $ws.Cells.Item($row, 3).Value = 5
}

$excel.CutCopyMode = 0

# Update the view state to match where the user ended up after adding data.
$ws.Activate()
$ws.Range("G103").Select()
$excel.ActiveWindow.ScrollRow = 84
